$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet logs one row per (date, quality) sample for the same market/product.
# A new week of sampling data (2022-02-04) is being inserted above the existing
# history, so push the existing rows 12:15 (and everything below them) down by
# four rows first.
$ws.Range("A12:T15").EntireRow.Insert()

# New sampling date for this week's entries.
$fecha = Get-Date -Year 2022 -Month 2 -Day 4 -Hour 0 -Minute 0 -Second 0

# Static values shared by every data row in this sheet.
$mercadoId = 1
$mercado   = "Agrícola del Norte S.A. de Arica"
$region    = "Arica y Parinacota"
$codreg    = 15
$tipo      = "Fruta"
$productoId = 100101
$producto   = "Berries"
$categoriaId = 100112025
$categoria   = "Frutilla"
$variedad    = "Sin especificar"
$unidad      = "`$/bandeja 3 kilos"
$origen      = "Región de Arica y Parinacota"
$kgUnidad    = 3

function Set-FilaFrutilla($Fila, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $PrecioKg) {
    $ws.Cells.Item($Fila, 1).Value = $mercadoId
    $ws.Cells.Item($Fila, 2).Value = $mercado
    $ws.Cells.Item($Fila, 3).Value = $region
    $ws.Cells.Item($Fila, 4).Value = $fecha
    $ws.Cells.Item($Fila, 5).Value = $codreg
    $ws.Cells.Item($Fila, 6).Value = $tipo
    $ws.Cells.Item($Fila, 7).Value = $productoId
    $ws.Cells.Item($Fila, 8).Value = $producto
    $ws.Cells.Item($Fila, 9).Value = $categoriaId
    $ws.Cells.Item($Fila, 10).Value = $categoria
    $ws.Cells.Item($Fila, 11).Value = $variedad
    $ws.Cells.Item($Fila, 12).Value = $Calidad
    $ws.Cells.Item($Fila, 13).Value = $Volumen
    $ws.Cells.Item($Fila, 14).Value = $PrecioMin
    $ws.Cells.Item($Fila, 15).Value = $PrecioMax
    $ws.Cells.Item($Fila, 16).Value = $PrecioProm
    $ws.Cells.Item($Fila, 17).Value = $unidad
    $ws.Cells.Item($Fila, 18).Value = $origen
    $ws.Cells.Item($Fila, 19).Value = $PrecioKg
    $ws.Cells.Item($Fila, 20).Value = $kgUnidad
}

Set-FilaFrutilla 12 "Especial" 100 8000 9000 8500 2833
Set-FilaFrutilla 13 "Primera"  130 6000 7000 6500 2167
Set-FilaFrutilla 14 "Segunda"  160 5000 6000 5500 1833
Set-FilaFrutilla 15 "Tercera"  100 4000 5000 4500 1500
